$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab07")

# --- Fix mis-decoded (mojibake) UTF-8 text in the Regional Economic Communities note (cell A103) ---
$ws.Range("A103").Value = "Regional Economic Communities: CEN-SAD = ""Community of Sahel-Saharan States""; COMESA = ""Common Market for Eastern and Southern Africa""; EAC = ""East African Community""; ECCAS = ""Economic Community of Central African States""; ECOWAS = ""Economic Community of West African States""; IGAD = ""Intergovernmental Authority on Development""; SADC = ""Southern African Development Community""; UMA = ""Arab Maghreb Union""; PALOP = ""Países Africanos de Língua Oficial Portuguesa""; ASEAN = ""Association of Southeast Asian Nations""; MERCOSUR = ""Mercado Común del Sur"". EU27 = ""European Union (27 members)"". OECD = ""Organisation for Economic Co-operation and Development""."

# --- Updated data values (refreshed ILOSTAT extraction) ---

$ws.Range("F6").Value = 89.141000000000005
$ws.Range("G6").Value = 89.9
$ws.Range("H6").Value = 88.548000000000002
$ws.Range("N6").Value = 65.998649999999998
$ws.Range("O6").Value = 15.199154
$ws.Range("P6").Value = 18.802070000000001
$ws.Range("F10").Value = 34.901000000000003
$ws.Range("G10").Value = 34.648000000000003
$ws.Range("H10").Value = 35.106999999999999
$ws.Range("N10").Value = 52.600037
$ws.Range("O10").Value = 23.536591000000001
$ws.Range("P10").Value = 23.863371999999998
$ws.Range("F13").Value = 74.333556000000002
$ws.Range("G13").Value = 76.146556000000004
$ws.Range("H13").Value = 72.531778000000003
$ws.Range("N13").Value = 67.038786000000002
$ws.Range("O13").Value = 13.785012
$ws.Range("P13").Value = 19.176203000000001
$ws.Range("F31").Value = 84.269000000000005
$ws.Range("G31").Value = 86.668000000000006
$ws.Range("H31").Value = 82.171000000000006
$ws.Range("N31").Value = 70.493742999999995
$ws.Range("O31").Value = 7.9363859999999997
$ws.Range("P31").Value = 21.569870999999999
$ws.Range("F32").Value = 14.222
$ws.Range("G32").Value = 5.84
$ws.Range("H32").Value = 23.016999999999999
$ws.Range("N32").Value = 83.971172999999993
$ws.Range("O32").Value = 1.6521459999999999
$ws.Range("P32").Value = 14.378831999999999
$ws.Range("F38").Value = 75.432417000000001
$ws.Range("G38").Value = 76.871082999999999
$ws.Range("H38").Value = 74.846999999999994
$ws.Range("N38").Value = 78.049736999999993
$ws.Range("O38").Value = 3.1143209999999999
$ws.Range("P38").Value = 18.835943
$ws.Range("F47").Value = 93.756
$ws.Range("G47").Value = 96.006
$ws.Range("H47").Value = 91.558999999999997
$ws.Range("N47").Value = 82.522205
$ws.Range("O47").Value = 2.1879460000000002
$ws.Range("P47").Value = 15.289866
$ws.Range("F61").Value = 89.778923000000006
$ws.Range("G61").Value = 92.056385000000006
$ws.Range("H61").Value = 87.559230999999997
$ws.Range("N61").Value = 82.919415999999998
$ws.Range("O61").Value = 2.5107970000000002
$ws.Range("P61").Value = 14.569784
$ws.Range("F62").Value = 80.322948999999994
$ws.Range("G62").Value = 81.687949000000003
$ws.Range("H62").Value = 78.998333000000002
$ws.Range("N62").Value = 74.190470000000005
$ws.Range("O62").Value = 4.8674289999999996
$ws.Range("P62").Value = 20.9421
$ws.Range("F63").Value = 37.622138
$ws.Range("G63").Value = 36.903137999999998
$ws.Range("H63").Value = 38.029915000000003
$ws.Range("N63").Value = 74.503746000000007
$ws.Range("O63").Value = 2.7750689999999998
$ws.Range("P63").Value = 22.728971999999999
$ws.Range("F64").Value = 52.890680000000003
$ws.Range("G64").Value = 50.63776
$ws.Range("H64").Value = 54.6614
$ws.Range("N64").Value = 76.343633999999994
$ws.Range("O64").Value = 3.9160789999999999
$ws.Range("P64").Value = 19.740286999999999
$ws.Range("F65").Value = 71.843050000000005
$ws.Range("G65").Value = 71.110399999999998
$ws.Range("H65").Value = 71.763999999999996
$ws.Range("N65").Value = 69.831230000000005
$ws.Range("O65").Value = 2.0036
$ws.Range("P65").Value = 28.165081000000001
$ws.Range("F66").Value = 50.143428999999998
$ws.Range("G66").Value = 50.035525999999997
$ws.Range("H66").Value = 50.043210999999999
$ws.Range("N66").Value = 74.455906999999996
$ws.Range("O66").Value = 3.0945830000000001
$ws.Range("P66").Value = 22.456108
$ws.Range("F67").Value = 73.047646999999998
$ws.Range("G67").Value = 73.670647000000002
$ws.Range("H67").Value = 72.421646999999993
$ws.Range("N67").Value = 71.180952000000005
$ws.Range("O67").Value = 3.1186859999999998
$ws.Range("P67").Value = 25.700364
$ws.Range("F68").Value = 85.493789000000007
$ws.Range("G68").Value = 87.206737000000004
$ws.Range("H68").Value = 83.850789000000006
$ws.Range("N68").Value = 74.210514000000003
$ws.Range("O68").Value = 2.965757
$ws.Range("P68").Value = 22.823726000000001
$ws.Range("F69").Value = 91.260999999999996
$ws.Range("G69").Value = 93.631600000000006
$ws.Range("H69").Value = 89.056600000000003
$ws.Range("N69").Value = 76.847256000000002
$ws.Range("O69").Value = 2.422129
$ws.Range("P69").Value = 20.730615
$ws.Range("F70").Value = 92.534999999999997
$ws.Range("G70").Value = 94.72775
$ws.Range("H70").Value = 90.491749999999996
$ws.Range("N70").Value = 76.413148000000007
$ws.Range("O70").Value = 3.3883350000000001
$ws.Range("P70").Value = 20.198516999999999
$ws.Range("F71").Value = 89.778923000000006
$ws.Range("G71").Value = 92.056385000000006
$ws.Range("H71").Value = 87.559230999999997
$ws.Range("N71").Value = 82.919415999999998
$ws.Range("O71").Value = 2.5107970000000002
$ws.Range("P71").Value = 14.569784
$ws.Range("F73").Value = 72.824332999999996
$ws.Range("G73").Value = 73.414867000000001
$ws.Range("H73").Value = 72.197599999999994
$ws.Range("N73").Value = 70.845065000000005
$ws.Range("O73").Value = 9.3947889999999994
$ws.Range("P73").Value = 19.760147
$ws.Range("F76").Value = 71.643570999999994
$ws.Range("G76").Value = 71.772285999999994
$ws.Range("H76").Value = 71.655000000000001
$ws.Range("N76").Value = 80.942971
$ws.Range("O76").Value = 1.310473
$ws.Range("P76").Value = 17.746555000000001
$ws.Range("F77").Value = 54.801614999999998
$ws.Range("G77").Value = 53.507154
$ws.Range("H77").Value = 55.623384999999999
$ws.Range("N77").Value = 76.992750000000001
$ws.Range("O77").Value = 4.6691349999999998
$ws.Range("P77").Value = 18.338115999999999
$ws.Range("F78").Value = 2.535593
$ws.Range("G78").Value = 2.8571110000000002
$ws.Range("H78").Value = 2.2556669999999999
$ws.Range("N78").Value = 82.303816999999995
$ws.Range("O78").Value = 4.6963670000000004
$ws.Range("P78").Value = 12.999814000000001
$ws.Range("F79").Value = 9.4875939999999996
$ws.Range("G79").Value = 10.129875
$ws.Range("H79").Value = 9.0340939999999996
$ws.Range("N79").Value = 79.953193999999996
$ws.Range("O79").Value = 3.7024180000000002
$ws.Range("P79").Value = 16.344387000000001
$ws.Range("F81").Value = 52.131875000000001
$ws.Range("G81").Value = 47.878749999999997
$ws.Range("H81").Value = 53.596874999999997
$ws.Range("N81").Value = 72.650583999999995
$ws.Range("O81").Value = 2.9982359999999999
$ws.Range("P81").Value = 24.351179999999999
$ws.Range("F82").Value = 79.723113999999995
$ws.Range("G82").Value = 81.023143000000005
$ws.Range("H82").Value = 78.427970999999999
$ws.Range("N82").Value = 74.102256999999994
$ws.Range("O82").Value = 4.7606830000000002
$ws.Range("P82").Value = 21.137060000000002
$ws.Range("F83").Value = 36.272395000000003
$ws.Range("G83").Value = 35.882151
$ws.Range("H83").Value = 36.581826
$ws.Range("N83").Value = 74.686588
$ws.Range("O83").Value = 2.75305
$ws.Range("P83").Value = 22.568918
$ws.Range("F84").Value = 92.718062000000003
$ws.Range("G84").Value = 94.816999999999993
$ws.Range("H84").Value = 90.662499999999994
$ws.Range("N84").Value = 78.403726000000006
$ws.Range("O84").Value = 2.2129240000000001
$ws.Range("P84").Value = 19.38335
$ws.Range("F86").Value = 78.876881999999995
$ws.Range("G86").Value = 80.922528999999997
$ws.Range("H86").Value = 76.953881999999993
$ws.Range("N86").Value = 74.629189999999994
$ws.Range("O86").Value = 3.8314859999999999
$ws.Range("P86").Value = 21.539321999999999
$ws.Range("F87").Value = 74.943888999999999
$ws.Range("G87").Value = 74.6905
$ws.Range("H87").Value = 74.855277999999998
$ws.Range("N87").Value = 69.686667
$ws.Range("O87").Value = 1.9551970000000001
$ws.Range("P87").Value = 28.358027
$ws.Range("F88").Value = 52.191249999999997
$ws.Range("G88").Value = 49.947249999999997
$ws.Range("H88").Value = 53.88
$ws.Range("N88").Value = 53.580717999999997
$ws.Range("O88").Value = 22.978207999999999
$ws.Range("P88").Value = 23.441074
$ws.Range("F89").Value = 49.622613000000001
$ws.Range("G89").Value = 47.557386999999999
$ws.Range("H89").Value = 50.778516000000003
$ws.Range("N89").Value = 73.166820000000001
$ws.Range("O89").Value = 3.4292769999999999
$ws.Range("P89").Value = 23.433306000000002
$ws.Range("F90").Value = 11.794364
$ws.Range("G90").Value = 11.468363999999999
$ws.Range("H90").Value = 12.101864000000001
$ws.Range("N90").Value = 83.264769000000001
$ws.Range("O90").Value = 3.3935559999999998
$ws.Range("P90").Value = 13.341673999999999
$ws.Range("F91").Value = 90.291269
$ws.Range("G91").Value = 92.975499999999997
$ws.Range("H91").Value = 88.131923
$ws.Range("N91").Value = 77.358963000000003
$ws.Range("O91").Value = 2.9361790000000001
$ws.Range("P91").Value = 19.704857000000001
$ws.Range("F93").Value = 59.106400000000001
$ws.Range("G93").Value = 57.348399999999998
$ws.Range("H93").Value = 60.719000000000001
$ws.Range("N93").Value = 74.189918000000006
$ws.Range("O93").Value = 3.0185080000000002
$ws.Range("P93").Value = 22.791574000000001
$ws.Range("F94").Value = 47.785842000000002
$ws.Range("G94").Value = 44.054420999999998
$ws.Range("H94").Value = 50.597842
$ws.Range("N94").Value = 77.660971000000004
$ws.Range("O94").Value = 2.629051
$ws.Range("P94").Value = 19.709983000000001
$ws.Range("F95").Value = 86.646083000000004
$ws.Range("G95").Value = 88.456333000000001
$ws.Range("H95").Value = 85.074250000000006
$ws.Range("N95").Value = 77.612714999999994
$ws.Range("O95").Value = 3.1140119999999998
$ws.Range("P95").Value = 19.273274000000001
$ws.Range("F96").Value = 62.347273000000001
$ws.Range("G96").Value = 62.116
$ws.Range("H96").Value = 62.626908999999998
$ws.Range("N96").Value = 70.418625000000006
$ws.Range("O96").Value = 2.8206669999999998
$ws.Range("P96").Value = 26.758692
$ws.Range("F97").Value = 88.406295999999998
$ws.Range("G97").Value = 91.305593000000002
$ws.Range("H97").Value = 86.029518999999993
$ws.Range("N97").Value = 79.867075
$ws.Range("O97").Value = 3.0619320000000001
$ws.Range("P97").Value = 17.070992
